# Before changing vmixclient to pointer
#
# 1) Add a new "Response" row to the Responses sheet (row 17):
#    id=99, type=9, Response.Text (same as other rows), new response text
#    "Is this working?"
# 2) Add a new "Shortcut" row to the Shortcuts sheet (row 8):
#    id=0, label "Merge Input=1"
# 3) Make the Responses sheet the active sheet/tab (it was Shortcuts before).

$wb = $excel.ActiveWorkbook

$responses = $wb.Worksheets.Item("Responses")
$responses.Cells.Item(17, 1).Value = 99
$responses.Cells.Item(17, 2).Value = 9
$responses.Cells.Item(17, 3).Value = "Response.Text"
$responses.Cells.Item(17, 4).Value = "Is this working?"

$shortcuts = $wb.Worksheets.Item("Shortcuts")
$shortcuts.Cells.Item(8, 1).Value = 0
$shortcuts.Cells.Item(8, 2).Value = "Merge Input=1"

# Shortcuts keeps its own last selection on row 8 even though it is no
# longer the active tab once we switch to Responses below.
$shortcuts.Activate()
$shortcuts.Range("B8").Select()

# Responses becomes the active sheet/tab (was Shortcuts before the edit),
# with the newly added row selected.
$responses.Activate()
$responses.Range("D17").Select()
